$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24 - this shifts the existing rows 24-41
# down to 25-42 (preserving their data/formatting) and leaves a fresh,
# empty row 24 (inheriting column D's date style from the row below).
$ws.Rows.Item(24).Insert()

# Populate the new row 24 with the new weekly price-report record.
$ws.Cells.Item(24, 1).Value = 11
$ws.Cells.Item(24, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(24, 3).Value = "Bíobío"
$ws.Cells.Item(24, 4).Value = 44484
$ws.Cells.Item(24, 5).Value = 8
$ws.Cells.Item(24, 6).Value = 100112001
$ws.Cells.Item(24, 7).Value = "Berenjena"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 230
$ws.Cells.Item(24, 11).Value = 5500
$ws.Cells.Item(24, 12).Value = 6000
$ws.Cells.Item(24, 13).Value = 5783
$ws.Cells.Item(24, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(24, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(24, 16).Value = 96
$ws.Cells.Item(24, 17).Value = 60
$ws.Cells.Item(24, 18).Value = "Hortaliza"
